$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("H8").Value = 71
$ws.Range("I8").Value = 64.666664
$ws.Range("J8").Value = 90
$ws.Range("K8").Value = 193.999992
$ws.Range("L8").Value = 270
$ws.Range("M8").Value = -54.99999199999999
$ws.Range("N8").Value = -548
$ws.Range("H9").Value = 56201.055
$ws.Range("I9").Value = 83481
$ws.Range("K9").Value = 83481
$ws.Range("M9").Value = -83312
$ws.Range("H32").Value = 2558.152
$ws.Range("J32").Value = 2774.1538
$ws.Range("L32").Value = 2774.1538
$ws.Range("N32").Value = -3426.1538
$ws.Range("H34").Value = 23750
$ws.Range("I34").Value = 23750
$ws.Range("K34").Value = 23750
$ws.Range("M34").Value = -23547
$ws.Range("H36").Value = 23750
$ws.Range("I36").Value = 23750
$ws.Range("K36").Value = 23750
$ws.Range("M36").Value = -23035
$ws.Range("H39").Value = 1910.5
$ws.Range("I39").Value = 298.14285
$ws.Range("J39").Value = 4167.8
$ws.Range("K39").Value = 894.4285500000001
$ws.Range("L39").Value = 12503.4
$ws.Range("M39").Value = -598.4285500000001
$ws.Range("N39").Value = -13095.4
$ws.Range("H41").Value = 685.5
$ws.Range("I41").Value = 324.5
$ws.Range("K41").Value = 324.5
$ws.Range("M41").Value = 115.5
$ws.Range("H43").Value = 2637.8
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2637.8
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2637.8
$ws.Range("M43").Value = $null
$ws.Range("N43").Value = -2775.8
$ws.Range("H69").Value = 16188
$ws.Range("I69").Value = 14519.25
$ws.Range("J69").Value = 17141.572
$ws.Range("K69").Value = 43557.75
$ws.Range("L69").Value = 51424.716
$ws.Range("M69").Value = -42683.75
$ws.Range("N69").Value = -53172.716
$ws.Range("H72").Value = 16188
$ws.Range("I72").Value = 14519.25
$ws.Range("J72").Value = 17141.572
$ws.Range("K72").Value = 130673.25
$ws.Range("L72").Value = 154274.148
$ws.Range("M72").Value = -126305.25
$ws.Range("N72").Value = -163010.148
$ws.Range("H88").Value = 2251.1667
$ws.Range("J88").Value = 2500.8
$ws.Range("L88").Value = 2500.8
$ws.Range("N88").Value = -3312.8
$ws.Range("H91").Value = 2251.1667
$ws.Range("J91").Value = 2500.8
$ws.Range("L91").Value = 2500.8
$ws.Range("N91").Value = -5308.8
$ws.Range("H98").Value = 49499.5
$ws.Range("J98").Value = 49499.5
$ws.Range("L98").Value = 49499.5
$ws.Range("N98").Value = -52495.5
$ws.Range("H100").Value = 3923.5
$ws.Range("I100").Value = 1597.125
$ws.Range("K100").Value = 1597.125
$ws.Range("M100").Value = -1056.125
$ws.Range("H101").Value = 4051.1052
$ws.Range("I101").Value = 1215.0834
$ws.Range("J101").Value = 8912.857
$ws.Range("K101").Value = 3645.2502
$ws.Range("L101").Value = 26738.571
$ws.Range("M101").Value = -2023.2502
$ws.Range("N101").Value = -29982.571
$ws.Range("H107").Value = 2622.44
$ws.Range("I107").Value = 1614.6364
$ws.Range("K107").Value = 1614.6364
$ws.Range("M107").Value = 305.3635999999999
$ws.Range("H116").Value = 19066.834
$ws.Range("J116").Value = 18186.318
$ws.Range("L116").Value = 18186.318
$ws.Range("N116").Value = -25070.318
$ws.Range("H122").Value = 49499.5
$ws.Range("J122").Value = 49499.5
$ws.Range("L122").Value = 148498.5
$ws.Range("N122").Value = -153398.5
$ws.Range("H132").Value = 3467.1052
$ws.Range("I132").Value = 3169.7742
$ws.Range("K132").Value = 9509.3226
$ws.Range("M132").Value = -6979.3226
$ws.Range("H135").Value = 1285.1818
$ws.Range("I135").Value = 1203.5238
$ws.Range("K135").Value = 10831.7142
$ws.Range("M135").Value = -8296.714199999999
$ws.Range("H137").Value = 2074.6875
$ws.Range("I137").Value = 1248.5714
$ws.Range("J137").Value = 2717.2222
$ws.Range("K137").Value = 3745.7142
$ws.Range("L137").Value = 8151.6666
$ws.Range("M137").Value = -1195.7142
$ws.Range("N137").Value = -13251.6666
$ws.Range("H138").Value = 3292.2449
$ws.Range("I138").Value = 3083.4546
$ws.Range("J138").Value = 3462.3704
$ws.Range("K138").Value = 9250.3638
$ws.Range("L138").Value = 10387.1112
$ws.Range("M138").Value = -4110.363799999999
$ws.Range("N138").Value = -20667.1112
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H32").Value = 2789.9592
$ws.Range("I32").Value = 2460.3333
$ws.Range("J32").Value = 6498.25
$ws.Range("K32").Value = 2460.3333
$ws.Range("L32").Value = 6498.25
$ws.Range("M32").Value = -2173.3333
$ws.Range("N32").Value = -7072.25
$ws.Range("H44").Value = 69987.336
$ws.Range("J44").Value = 69987.336
$ws.Range("L44").Value = 69987.336
$ws.Range("N44").Value = -70963.336
$ws.Range("H45").Value = 1919.1177
$ws.Range("I45").Value = 1364
$ws.Range("K45").Value = 1364
$ws.Range("M45").Value = -987
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("H61").Value = 14363.053
$ws.Range("I61").Value = 14637
$ws.Range("K61").Value = 14637
$ws.Range("M61").Value = -14425
$ws.Range("H74").Value = 4542.5293
$ws.Range("I74").Value = 4105.7144
$ws.Range("J74").Value = 6581
$ws.Range("K74").Value = 4105.7144
$ws.Range("L74").Value = 6581
$ws.Range("M74").Value = -3231.7144
$ws.Range("N74").Value = -8329
$ws.Range("H77").Value = 4542.5293
$ws.Range("I77").Value = 4105.7144
$ws.Range("J77").Value = 6581
$ws.Range("K77").Value = 20528.572
$ws.Range("L77").Value = 32905
$ws.Range("M77").Value = -16160.572
$ws.Range("N77").Value = -41641
$ws.Range("H86").Value = 24000
$ws.Range("I86").Value = 24000
$ws.Range("K86").Value = 24000
$ws.Range("M86").Value = -22814
$ws.Range("H88").Value = 4759.7095
$ws.Range("I88").Value = 1080.8
$ws.Range("J88").Value = 6511.5713
$ws.Range("K88").Value = 1080.8
$ws.Range("L88").Value = 6511.5713
$ws.Range("M88").Value = -674.8
$ws.Range("N88").Value = -7323.5713
$ws.Range("H89").Value = 24000
$ws.Range("I89").Value = 24000
$ws.Range("K89").Value = 72000
$ws.Range("M89").Value = -66072
$ws.Range("H91").Value = 4759.7095
$ws.Range("I91").Value = 1080.8
$ws.Range("J91").Value = 6511.5713
$ws.Range("K91").Value = 1080.8
$ws.Range("L91").Value = 6511.5713
$ws.Range("M91").Value = 323.2
$ws.Range("N91").Value = -9319.5713
$ws.Range("H97").Value = 11649
$ws.Range("I97").Value = 15594.857
$ws.Range("J97").Value = 6124.8
$ws.Range("K97").Value = 15594.857
$ws.Range("L97").Value = 6124.8
$ws.Range("M97").Value = -15098.857
$ws.Range("N97").Value = -7116.8
$ws.Range("H102").Value = 1992.9231
$ws.Range("I102").Value = 1992.9231
$ws.Range("K102").Value = 1992.9231
$ws.Range("M102").Value = -370.9231
$ws.Range("H110").Value = 1824.9
$ws.Range("J110").Value = 2749
$ws.Range("L110").Value = 2749
$ws.Range("N110").Value = -6839
$ws.Range("H132").Value = 3386.25
$ws.Range("I132").Value = 2216.5908
$ws.Range("J132").Value = 7675
$ws.Range("K132").Value = 6649.7724
$ws.Range("L132").Value = 23025
$ws.Range("M132").Value = -4119.7724
$ws.Range("N132").Value = -28085
$ws.Range("H136").Value = 14363.053
$ws.Range("I136").Value = 14637
$ws.Range("K136").Value = 43911
$ws.Range("M136").Value = -41361

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 1821.52
$ws.Range("I20").Value = 1387.5555
$ws.Range("J20").Value = 2937.4285
$ws.Range("K20").Value = 1387.5555
$ws.Range("L20").Value = 2937.4285
$ws.Range("M20").Value = -1140.5555
$ws.Range("N20").Value = -3431.4285
$ws.Range("H21").Value = 18588
$ws.Range("J21").Value = 18588
$ws.Range("L21").Value = 18588
$ws.Range("N21").Value = -19060
$ws.Range("H92").Value = 991
$ws.Range("J92").Value = 991
$ws.Range("L92").Value = 991
$ws.Range("N92").Value = -5983
$ws.Range("H107").Value = 1619.8667
$ws.Range("I107").Value = 1671.5264
$ws.Range("K107").Value = 1671.5264
$ws.Range("M107").Value = 248.4736
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H134").Value = 11052.36
$ws.Range("I134").Value = 12313.909
$ws.Range("J134").Value = 10061.143
$ws.Range("K134").Value = 36941.727
$ws.Range("L134").Value = 30183.429
$ws.Range("M134").Value = -34406.727
$ws.Range("N134").Value = -35253.429

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2977.4614
$ws.Range("I31").Value = 1199.4
$ws.Range("K31").Value = 1199.4
$ws.Range("M31").Value = -904.4000000000001
$ws.Range("H34").Value = 2977.4614
$ws.Range("I34").Value = 1199.4
$ws.Range("K34").Value = 1199.4
$ws.Range("M34").Value = -997.4000000000001
$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -26250
$ws.Range("H51").Value = 54999
$ws.Range("J51").Value = 54999
$ws.Range("L51").Value = 54999
$ws.Range("N51").Value = -56471
$ws.Range("H58").Value = 5179.0557
$ws.Range("I58").Value = 4291.8184
$ws.Range("J58").Value = 6573.2856
$ws.Range("K58").Value = 4291.8184
$ws.Range("L58").Value = 6573.2856
$ws.Range("M58").Value = -4088.8184
$ws.Range("N58").Value = -6979.2856
$ws.Range("H59").Value = 33398.734
$ws.Range("I59").Value = 8499.333
$ws.Range("J59").Value = 49998.332
$ws.Range("K59").Value = 8499.333
$ws.Range("L59").Value = 49998.332
$ws.Range("M59").Value = -7354.333000000001
$ws.Range("N59").Value = -52288.332
$ws.Range("H60").Value = 31515.807
$ws.Range("I60").Value = 13400
$ws.Range("J60").Value = 34999.617
$ws.Range("K60").Value = 13400
$ws.Range("L60").Value = 34999.617
$ws.Range("M60").Value = -12889
$ws.Range("N60").Value = -36021.617
$ws.Range("H61").Value = 54999
$ws.Range("J61").Value = 54999
$ws.Range("L61").Value = 54999
$ws.Range("N61").Value = -55695
$ws.Range("H62").Value = 18114.875
$ws.Range("I62").Value = 5509.6665
$ws.Range("J62").Value = 25678
$ws.Range("K62").Value = 5509.6665
$ws.Range("L62").Value = 25678
$ws.Range("M62").Value = -4885.6665
$ws.Range("N62").Value = -26926
$ws.Range("H65").Value = 18114.875
$ws.Range("I65").Value = 5509.6665
$ws.Range("J65").Value = 25678
$ws.Range("K65").Value = 27548.3325
$ws.Range("L65").Value = 128390
$ws.Range("M65").Value = -24428.3325
$ws.Range("N65").Value = -134630
$ws.Range("H96").Value = 9868.667
$ws.Range("J96").Value = 9868.667
$ws.Range("L96").Value = 9868.667
$ws.Range("N96").Value = -15360.667
$ws.Range("H99").Value = 2988.923
$ws.Range("I99").Value = 2921.3125
$ws.Range("J99").Value = 3097.1
$ws.Range("K99").Value = 2921.3125
$ws.Range("L99").Value = 3097.1
$ws.Range("M99").Value = -1423.3125
$ws.Range("N99").Value = -6093.1
$ws.Range("H122").Value = 5130.4116
$ws.Range("J122").Value = 5202.4
$ws.Range("L122").Value = 15607.2
$ws.Range("N122").Value = -20507.2
$ws.Range("H126").Value = 2988.923
$ws.Range("I126").Value = 2921.3125
$ws.Range("J126").Value = 3097.1
$ws.Range("K126").Value = 8763.9375
$ws.Range("L126").Value = 9291.3
$ws.Range("M126").Value = -6293.9375
$ws.Range("N126").Value = -14231.3
$ws.Range("H132").Value = 2023.96
$ws.Range("I132").Value = 2056.2083
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 6168.624899999999
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = -3638.624899999999
$ws.Range("N132").Value = -8810
$ws.Range("H134").Value = 5530.585
$ws.Range("I134").Value = 5077.2563
$ws.Range("J134").Value = 6793.4287
$ws.Range("K134").Value = 15231.7689
$ws.Range("L134").Value = 20380.2861
$ws.Range("M134").Value = -12696.7689
$ws.Range("N134").Value = -25450.2861
$ws.Range("H136").Value = 5179.0557
$ws.Range("I136").Value = 4291.8184
$ws.Range("J136").Value = 6573.2856
$ws.Range("K136").Value = 12875.4552
$ws.Range("L136").Value = 19719.8568
$ws.Range("M136").Value = -10325.4552
$ws.Range("N136").Value = -24819.8568

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H6").Value = 247.5
$ws.Range("I6").Value = 295.66666
$ws.Range("K6").Value = 886.9999799999999
$ws.Range("M6").Value = -773.9999799999999
$ws.Range("H14").Value = 404.2
$ws.Range("I14").Value = 404.2
$ws.Range("K14").Value = 1212.6
$ws.Range("M14").Value = -1039.6
$ws.Range("H40").Value = 274.2143
$ws.Range("I40").Value = 78.333336
$ws.Range("J40").Value = 421.125
$ws.Range("K40").Value = 313.333344
$ws.Range("L40").Value = 1684.5
$ws.Range("M40").Value = -244.333344
$ws.Range("N40").Value = -1822.5
$ws.Range("H62").Value = 8713.714
$ws.Range("J62").Value = 9699.2
$ws.Range("L62").Value = 29097.6
$ws.Range("N62").Value = -30469.6
$ws.Range("H65").Value = 8713.714
$ws.Range("J65").Value = 9699.2
$ws.Range("L65").Value = 87292.8
$ws.Range("N65").Value = -94156.8
$ws.Range("H74").Value = 10499
$ws.Range("J74").Value = 10499
$ws.Range("L74").Value = 31497
$ws.Range("N74").Value = -33619
$ws.Range("H77").Value = 10499
$ws.Range("J77").Value = 10499
$ws.Range("L77").Value = 94491
$ws.Range("N77").Value = -105099
$ws.Range("H98").Value = 38465840
$ws.Range("I98").Value = 3417
$ws.Range("J98").Value = 55560252
$ws.Range("K98").Value = 10251
$ws.Range("L98").Value = 166680756
$ws.Range("M98").Value = -8753
$ws.Range("N98").Value = -166683752
$ws.Range("H129").Value = 12823054
$ws.Range("I129").Value = 1354.5
$ws.Range("J129").Value = 23813082
$ws.Range("K129").Value = 4063.5
$ws.Range("L129").Value = 71439246
$ws.Range("M129").Value = 936.5
$ws.Range("N129").Value = -71449246

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 124.2
$ws.Range("I2").Value = 84.55
$ws.Range("J2").Value = 282.8
$ws.Range("K2").Value = 84.55
$ws.Range("L2").Value = 282.8
$ws.Range("M2").Value = 28.45
$ws.Range("N2").Value = -508.8
$ws.Range("H7").Value = 950001
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 950001
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 950001
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = -950225
$ws.Range("H8").Value = 950001
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 950001
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 950001
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = -950279
$ws.Range("H102").Value = 7189.8184
$ws.Range("I102").Value = 6875.625
$ws.Range("J102").Value = 8027.6665
$ws.Range("K102").Value = 6875.625
$ws.Range("L102").Value = 8027.6665
$ws.Range("M102").Value = -5253.625
$ws.Range("N102").Value = -11271.6665
$ws.Range("H126").Value = 6418.1816
$ws.Range("I126").Value = 3100
$ws.Range("K126").Value = 9300
$ws.Range("M126").Value = -6830
$ws.Range("H132").Value = 5168.8667
$ws.Range("I132").Value = 3338.5557
$ws.Range("K132").Value = 10015.6671
$ws.Range("M132").Value = -7485.667099999999

# Sheet index 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 465.66666
$ws.Range("I16").Value = 449.5
$ws.Range("K16").Value = 449.5
$ws.Range("M16").Value = -279.5
$ws.Range("H22").Value = 1040.7059
$ws.Range("I22").Value = 679.5
$ws.Range("J22").Value = 1088.8667
$ws.Range("K22").Value = 679.5
$ws.Range("L22").Value = 1088.8667
$ws.Range("M22").Value = -384.5
$ws.Range("N22").Value = -1678.8667
$ws.Range("H27").Value = 1040.7059
$ws.Range("I27").Value = 679.5
$ws.Range("J27").Value = 1088.8667
$ws.Range("K27").Value = 679.5
$ws.Range("L27").Value = 1088.8667
$ws.Range("M27").Value = -572.5
$ws.Range("N27").Value = -1302.8667
$ws.Range("H46").Value = 2162.3462
$ws.Range("J46").Value = 2551.2
$ws.Range("L46").Value = 2551.2
$ws.Range("N46").Value = -2927.2
$ws.Range("H55").Value = 471.7857
$ws.Range("I55").Value = 151.625
$ws.Range("J55").Value = 898.6667
$ws.Range("K55").Value = 151.625
$ws.Range("L55").Value = 898.6667
$ws.Range("M55").Value = 21.375
$ws.Range("N55").Value = -1244.6667
$ws.Range("H61").Value = 2700.4285
$ws.Range("I61").Value = 2700.4285
$ws.Range("K61").Value = 2700.4285
$ws.Range("M61").Value = -2498.4285
$ws.Range("H82").Value = 5910.222
$ws.Range("I82").Value = 1676
$ws.Range("K82").Value = 1676
$ws.Range("M82").Value = -1315
$ws.Range("H85").Value = 5910.222
$ws.Range("I85").Value = 1676
$ws.Range("K85").Value = 1676
$ws.Range("M85").Value = -428
$ws.Range("H93").Value = 1527.8
$ws.Range("I93").Value = 1231.6666
$ws.Range("K93").Value = 1231.6666
$ws.Range("M93").Value = 16.33339999999998
$ws.Range("H113").Value = 2700.4285
$ws.Range("I113").Value = 2700.4285
$ws.Range("K113").Value = 2700.4285
$ws.Range("M113").Value = -530.4285
$ws.Range("H122").Value = 5836.1377
$ws.Range("I122").Value = 3931.875
$ws.Range("J122").Value = 6561.5713
$ws.Range("K122").Value = 11795.625
$ws.Range("L122").Value = 19684.7139
$ws.Range("M122").Value = -9345.625
$ws.Range("N122").Value = -24584.7139
$ws.Range("H132").Value = 4312.6665
$ws.Range("I132").Value = 3337.3333
$ws.Range("K132").Value = 10011.9999
$ws.Range("M132").Value = -7481.999899999999
$ws.Range("H136").Value = 1566.6
$ws.Range("I136").Value = 1566.6
$ws.Range("K136").Value = 4699.799999999999
$ws.Range("M136").Value = -2149.799999999999

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 9956.333
$ws.Range("I62").Value = 9970
$ws.Range("J62").Value = 9949.5
$ws.Range("K62").Value = 9970
$ws.Range("L62").Value = 9949.5
$ws.Range("M62").Value = -9346
$ws.Range("N62").Value = -11197.5
$ws.Range("H65").Value = 9956.333
$ws.Range("I65").Value = 9970
$ws.Range("J65").Value = 9949.5
$ws.Range("K65").Value = 49850
$ws.Range("L65").Value = 49747.5
$ws.Range("M65").Value = -46730
$ws.Range("N65").Value = -55987.5
$ws.Range("H81").Value = 4027.3547
$ws.Range("I81").Value = 4384.0454
$ws.Range("J81").Value = 3155.4443
$ws.Range("K81").Value = 8768.0908
$ws.Range("L81").Value = 6310.8886
$ws.Range("M81").Value = -7707.0908
$ws.Range("N81").Value = -8432.8886
$ws.Range("H84").Value = 4027.3547
$ws.Range("I84").Value = 4384.0454
$ws.Range("J84").Value = 3155.4443
$ws.Range("K84").Value = 43840.454
$ws.Range("L84").Value = 31554.443
$ws.Range("M84").Value = -38536.454
$ws.Range("N84").Value = -42162.443
$ws.Range("H104").Value = 18456.334
$ws.Range("J104").Value = 18456.334
$ws.Range("L104").Value = 18456.334
$ws.Range("N104").Value = -25444.334
$ws.Range("H107").Value = 775.1739
$ws.Range("I107").Value = 753.2222
$ws.Range("J107").Value = 854.2
$ws.Range("K107").Value = 2259.6666
$ws.Range("L107").Value = 2562.6
$ws.Range("M107").Value = -339.6666
$ws.Range("N107").Value = -6402.6
$ws.Range("H113").Value = 252
$ws.Range("I113").Value = 252
$ws.Range("K113").Value = 756
$ws.Range("M113").Value = 1414
$ws.Range("H132").Value = 4932.241
$ws.Range("I132").Value = 4912.6816
$ws.Range("J132").Value = 4993.7144
$ws.Range("K132").Value = 14738.0448
$ws.Range("L132").Value = 14981.1432
$ws.Range("M132").Value = -12208.0448
$ws.Range("N132").Value = -20041.1432
$ws.Range("H136").Value = 17756.5
$ws.Range("I136").Value = 17756.5
$ws.Range("K136").Value = 53269.5
$ws.Range("M136").Value = -50719.5

